# Refresh the cryptos list with the latest coinranking.com snapshot:
# updates the Price (D) and Volume(1h) (E) columns, and corrects the
# swapped BinanceUSD / LidoDAOToken rows (37-38).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.465.19"
$ws.Range("E2").Value = "  -2.56%  "
$ws.Range("D3").Value = "1.983.71"
$ws.Range("E3").Value = "  -3.24%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'244.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.24%  "
$ws.Range("D6").Value = "'0.627"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.72%  "
$ws.Range("D7").Value = "'58.86"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -9.93%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.375"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").Value = "'57.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.27%  "
$ws.Range("E11").Value = "  +6.04%  "
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").Value = "'23.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.52%  "
$ws.Range("D14").Value = "'0.862"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.71%  "
$ws.Range("D15").Value = "'14.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.97%  "
$ws.Range("D16").Value = "2.275.79"
$ws.Range("E16").Value = "  -3.32%  "
$ws.Range("D17").Value = "'5.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.15%  "
$ws.Range("D18").Value = "1.988.02"
$ws.Range("E18").Value = "  -3.11%  "
$ws.Range("D19").Value = "36.373.23"
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("D20").Value = "'70.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.15%  "
$ws.Range("D21").Value = "0.0₃0864"
$ws.Range("E21").Value = "  -1.52%  "
$ws.Range("D22").Value = "'5.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.57%  "
$ws.Range("D23").Value = "'234.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("D26").Value = "'2.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.98%  "
$ws.Range("D27").Value = "'10.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.31%  "
$ws.Range("D28").Value = "'162.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "'19.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.62%  "
$ws.Range("D30").Value = "'0.130"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.35%  "
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("D33").Value = "'4.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.25%  "
$ws.Range("D34").Value = "'0.0628"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("D35").Value = "'4.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.10%  "
$ws.Range("D36").Value = "'6.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.71%  "
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.40%  "
$ws.Range("E39").Value = "  -3.80%  "
$ws.Range("D40").Value = "'3.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.87%  "
$ws.Range("D41").Value = "'1.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("D42").Value = "'0.0958"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.69%  "
$ws.Range("E43").Value = "  -3.65%  "
$ws.Range("D44").Value = "'0.0213"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("E45").Value = "  -4.80%  "
$ws.Range("D46").Value = "'92.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.08%  "
$ws.Range("D47").Value = "'16.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.94%  "
$ws.Range("D48").Value = "1.385.94"
$ws.Range("E48").Value = "  -2.49%  "
$ws.Range("D49").Value = "'7.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.50%  "
$ws.Range("D50").Value = "'2.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.04%  "
$ws.Range("D51").Value = "'45.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.40%  "

Write-Output "Applied 93 cell updates"
